# Generate Report for Handoff
#
# The two separate source files (affb26ab-...md and f8878dfd-...md) were
# consolidated into a single handoff package (3d6d923a-...md /
# ffff92659b39-...md) that is now "Ready for handoff". The per-language
# sheets' "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns (E/F) are dropped since nothing has been handed back
# yet, and the Latest Handoff File / Latest Handoff Datetime columns now
# both point at the new consolidated handoff package.

$wb = $excel.ActiveWorkbook

$oldFile1 = "affb26ab-e351-47c8-b895-683175176dd9.md"
$oldFile2 = "f8878dfd-f9fe-4b6f-8dde-91c59f7ccd04.md"
$newFile1 = "3d6d923a-4f6f-4169-992f-ccc384019ff3.md"
$newFile2 = "ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md"

$newStatus = "Ready for handoff"

function Set-HyperlinkDisplay($ws, $cellRef, $address, $display) {
    # Remove any existing hyperlink anchored at this cell (searching by
    # matching range address), preserving the others, then re-add one
    # pointing at $address with the requested display text.
    $target = $ws.Range($cellRef)
    $toDelete = @()
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $target.Address()) {
            $toDelete += $h
        }
    }
    foreach ($h in $toDelete) {
        $h.Delete()
    }
    $ws.Hyperlinks.Add($target, $address, "", "", $display) | Out-Null
}

# ---------------------------------------------------------------------
# Overview sheet: only the hyperlink display text (and, through the
# shared strings that the other sheets reference, the Status text) needs
# to move on to the new filenames.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-HyperlinkDisplay $wsOverview "A2" "https://github.com/OpenLocalizationTest/oltest/blob/19438f37ce5bd052ea6e7091bb3fdd8a095b7e52/e2e/$newFile1" $newFile1
Set-HyperlinkDisplay $wsOverview "A3" "https://github.com/OpenLocalizationTest/oltest/blob/19438f37ce5bd052ea6e7091bb3fdd8a095b7e52/e2e/$newFile2" $newFile2

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): update source/handoff file names,
# handoff datetime, status, and drop the "Latest Target File" /
# "Latest Handback File" columns (E/F) since nothing has been handed
# back for the new package yet.
# ---------------------------------------------------------------------
function Update-LangSheet($sheetName, $xlfLang, $handoffDatetimeRow2Row3, $handoffUrlPrefix) {
    $ws = $wb.Worksheets.Item($sheetName)

    $newXlf = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.$xlfLang.xlf"

    Set-HyperlinkDisplay $ws "A2" "https://github.com/OpenLocalizationTest/oltest/blob/19438f37ce5bd052ea6e7091bb3fdd8a095b7e52/e2e/$newFile1" $newFile1
    Set-HyperlinkDisplay $ws "C2" "$handoffUrlPrefix/$newXlf" $newXlf
    Set-HyperlinkDisplay $ws "A3" "https://github.com/OpenLocalizationTest/oltest/blob/19438f37ce5bd052ea6e7091bb3fdd8a095b7e52/e2e/$newFile2" $newFile2
    Set-HyperlinkDisplay $ws "C3" "$handoffUrlPrefix/$newXlf" $newXlf

    $ws.Range("B2").Value = $newStatus
    $ws.Range("D2").Value = $handoffDatetimeRow2Row3
    $ws.Range("E2").ClearContents()
    $ws.Range("F2").ClearContents()
    $ws.Range("G2").Value = "0001-01-01 00:00:00"

    $ws.Range("B3").Value = $newStatus
    $ws.Range("D3").Value = $handoffDatetimeRow2Row3
    $ws.Range("E3").ClearContents()
    $ws.Range("F3").ClearContents()
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
}

Update-LangSheet "zh-cn" "zh-cn" "2016-03-09 14:29:08" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9941e9bb8ebec7c6d2dc21b6403520cca56624d8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
Update-LangSheet "de-de" "de-de" "2016-03-09 14:29:11" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97ffd96e7582217760e789e4d051fa10b886dc96/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"
